$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44259
$ws.Range("J2").Value = 30

$ws.Range("D3").Value = 44313
$ws.Range("J3").Value = 20

$ws.Range("D4").Value = 44176
$ws.Range("J4").Value = 10

$ws.Range("D5").Value = 44365
$ws.Range("J5").Value = 55
$ws.Range("K5").Value = 5000
$ws.Range("L5").Value = 5000
$ws.Range("M5").Value = 5000
$ws.Range("P5").Value = 5000

$ws.Range("D6").Value = 44315
$ws.Range("J6").Value = 40

$ws.Range("D7").Value = 44316
$ws.Range("J7").Value = 20
$ws.Range("K7").Value = 4000
$ws.Range("L7").Value = 4000
$ws.Range("M7").Value = 4000
$ws.Range("P7").Value = 4000

$ws.Range("D8").Value = 44291
$ws.Range("J8").Value = 35
$ws.Range("K8").Value = 4000
$ws.Range("L8").Value = 4000
$ws.Range("M8").Value = 4000
$ws.Range("P8").Value = 4000

$ws.Range("D9").Value = 44280
$ws.Range("J9").Value = 55

$ws.Range("D10").Value = 44301
$ws.Range("J10").Value = 40
$ws.Range("K10").Value = 3000
$ws.Range("L10").Value = 3000
$ws.Range("M10").Value = 3000
$ws.Range("P10").Value = 3000

$ws.Range("D11").Value = 44390
$ws.Range("K11").Value = 6000
$ws.Range("L11").Value = 6000
$ws.Range("M11").Value = 6000
$ws.Range("P11").Value = 6000

$ws.Range("D12").Value = 44312
$ws.Range("J12").Value = 50
